# new Madigan bike hours
# Update Riders (C) and Average (D) values for the weekly Ridership sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Row 2
$ws.Range("C2").Value = 238
$ws.Range("D2").Value = 212.17

# Row 3
$ws.Range("C3").Value = 177
$ws.Range("D3").Value = 194.17

# Row 4
$ws.Range("C4").Value = 275
$ws.Range("D4").Value = 225.09

# Row 5
$ws.Range("C5").Value = 246
$ws.Range("D5").Value = 237.33

# Row 6
$ws.Range("C6").Value = 176
$ws.Range("D6").Value = 121.92

# Row 7
$ws.Range("C7").Value = 125
$ws.Range("D7").Value = 103.64
